# Generate Report for Handback
# - Updates the "Status" column text on all three sheets from
#   "Ready for handoff" to "Handed back: in sync with en-US"
# - Populates "Latest Target File" (E) and "Latest Handback File" (F)
#   for the two tracked source files on the zh-cn and de-de sheets,
#   each mirroring the existing Source File / Handoff File hyperlink.
# - Stamps "Latest Handback DateTime" (G) with the handback timestamp
#   for each language.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B2").Value = $newStatus
$zh.Range("B3").Value = $newStatus

$de = $wb.Worksheets.Item("de-de")
$de.Range("B2").Value = $newStatus
$de.Range("B3").Value = $newStatus

# Cornflower blue (FF6495ED), packed BGR for Font.Color
$hyperlinkColor = 15570276

function Style-AsHyperlink($rng) {
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Underline = $true
    $rng.Font.Color = $hyperlinkColor
}

# ---- 2. zh-cn sheet: target file / handback file / handback datetime ----
$zh.Range("E2").Value = "52ed64bd-7ad3-4b35-98a5-bec23f741d3f.md"
$null = $zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/ec08001f30795f72e2abb9d1b5138c4b7c8097de/e2e/52ed64bd-7ad3-4b35-98a5-bec23f741d3f.md", "", "", "52ed64bd-7ad3-4b35-98a5-bec23f741d3f.md")
Style-AsHyperlink $zh.Range("E2")

$zh.Range("F2").Value = "52ed64bd-7ad3-4b35-98a5-bec23f741d3f.82dbc4e45d244c08ba4134c49868731c4d30bd08.zh-cn.xlf"
$null = $zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a2ff3db881634b53e7e4f093733372b1a7ceed86/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/52ed64bd-7ad3-4b35-98a5-bec23f741d3f.82dbc4e45d244c08ba4134c49868731c4d30bd08.zh-cn.xlf", "", "", "52ed64bd-7ad3-4b35-98a5-bec23f741d3f.82dbc4e45d244c08ba4134c49868731c4d30bd08.zh-cn.xlf")
Style-AsHyperlink $zh.Range("F2")

$zh.Range("G2").Value = "2016-01-27 08:33:18"

$zh.Range("E3").Value = "5db9de63-b467-4b11-8212-54a0cf1a683f.md"
$null = $zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/ec08001f30795f72e2abb9d1b5138c4b7c8097de/e2e/5db9de63-b467-4b11-8212-54a0cf1a683f.md", "", "", "5db9de63-b467-4b11-8212-54a0cf1a683f.md")
Style-AsHyperlink $zh.Range("E3")

$zh.Range("F3").Value = "5db9de63-b467-4b11-8212-54a0cf1a683f.06dfc7c78e85c7a2245ee7fa47884b41bd3bd518.zh-cn.xlf"
$null = $zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a2ff3db881634b53e7e4f093733372b1a7ceed86/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/5db9de63-b467-4b11-8212-54a0cf1a683f.06dfc7c78e85c7a2245ee7fa47884b41bd3bd518.zh-cn.xlf", "", "", "5db9de63-b467-4b11-8212-54a0cf1a683f.06dfc7c78e85c7a2245ee7fa47884b41bd3bd518.zh-cn.xlf")
Style-AsHyperlink $zh.Range("F3")

$zh.Range("G3").Value = "2016-01-27 08:33:18"

# ---- 3. de-de sheet: target file / handback file / handback datetime ----
$de.Range("E2").Value = "52ed64bd-7ad3-4b35-98a5-bec23f741d3f.md"
$null = $de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/ec08001f30795f72e2abb9d1b5138c4b7c8097de/e2e/52ed64bd-7ad3-4b35-98a5-bec23f741d3f.md", "", "", "52ed64bd-7ad3-4b35-98a5-bec23f741d3f.md")
Style-AsHyperlink $de.Range("E2")

$de.Range("F2").Value = "52ed64bd-7ad3-4b35-98a5-bec23f741d3f.82dbc4e45d244c08ba4134c49868731c4d30bd08.de-de.xlf"
$null = $de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5054d915043b89e66fd03f479dbe5d8358f5f55/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/52ed64bd-7ad3-4b35-98a5-bec23f741d3f.82dbc4e45d244c08ba4134c49868731c4d30bd08.de-de.xlf", "", "", "52ed64bd-7ad3-4b35-98a5-bec23f741d3f.82dbc4e45d244c08ba4134c49868731c4d30bd08.de-de.xlf")
Style-AsHyperlink $de.Range("F2")

$de.Range("G2").Value = "2016-01-27 08:33:43"

$de.Range("E3").Value = "5db9de63-b467-4b11-8212-54a0cf1a683f.md"
$null = $de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/ec08001f30795f72e2abb9d1b5138c4b7c8097de/e2e/5db9de63-b467-4b11-8212-54a0cf1a683f.md", "", "", "5db9de63-b467-4b11-8212-54a0cf1a683f.md")
Style-AsHyperlink $de.Range("E3")

$de.Range("F3").Value = "5db9de63-b467-4b11-8212-54a0cf1a683f.06dfc7c78e85c7a2245ee7fa47884b41bd3bd518.de-de.xlf"
$null = $de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5054d915043b89e66fd03f479dbe5d8358f5f55/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/5db9de63-b467-4b11-8212-54a0cf1a683f.06dfc7c78e85c7a2245ee7fa47884b41bd3bd518.de-de.xlf", "", "", "5db9de63-b467-4b11-8212-54a0cf1a683f.06dfc7c78e85c7a2245ee7fa47884b41bd3bd518.de-de.xlf")
Style-AsHyperlink $de.Range("F3")

$de.Range("G3").Value = "2016-01-27 08:33:43"

Write-Host "Handback report applied."
